$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.611.60'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '3.446.46'
$ws.Range('E3').Value = '  -3.61%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.47'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.77'
$ws.Range('E6').Value = '  -7.34%  '
$ws.Range('D7').Value = '3.445.55'
$ws.Range('E7').Value = '  -3.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.47'
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('E11').Value = '  -9.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.380'
$ws.Range('E12').Value = '  -7.48%  '
$ws.Range('D13').Value = '4.033.03'
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('E14').Value = '  -10.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.64'
$ws.Range('E15').Value = '  -9.57%  '
$ws.Range('D16').Value = '3.445.60'
$ws.Range('E16').Value = '  -3.77%  '
$ws.Range('D17').Value = '65.584.78'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.94'
$ws.Range('E19').Value = '  -9.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.78'
$ws.Range('E20').Value = '  -8.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.75'
$ws.Range('E21').Value = '  -6.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '395.42'
$ws.Range('E22').Value = '  -6.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.549'
$ws.Range('E23').Value = '  -9.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.36'
$ws.Range('E24').Value = '  -6.14%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '3.593.28'
$ws.Range('E26').Value = '  -3.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000106'
$ws.Range('E27').Value = '  -9.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  -8.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.18'
$ws.Range('E30').Value = '  -10.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.22'
$ws.Range('E31').Value = '  -11.73%  '
$ws.Range('D32').Value = '3.454.96'
$ws.Range('E32').Value = '  -3.49%  '
$ws.Range('E34').Value = '  -6.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '22.93'
$ws.Range('E35').Value = '  -7.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '172.88'
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -13.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.92'
$ws.Range('E38').Value = '  -10.14%  '
$ws.Range('E39').Value = '  -7.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.80'
$ws.Range('E40').Value = '  -12.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0778'
$ws.Range('E41').Value = '  -8.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.820'
$ws.Range('E42').Value = '  -6.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.46'
$ws.Range('E43').Value = '  -5.52%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.45'
$ws.Range('E45').Value = '  -13.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.63'
$ws.Range('E46').Value = '  -11.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.17'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.10'
$ws.Range('E48').Value = '  -2.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.57'
$ws.Range('E49').Value = '  -7.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.10'
$ws.Range('E50').Value = '  -15.57%  '
$ws.Range('D51').Value = '2.219.90'
$ws.Range('E51').Value = '  -7.47%  '
